# Fruta / hortaliza, semanal
# Re-shuffle the data rows (2-26) across columns D, I, J, K, L, M, O, P
# following the permutation observed between the old and new workbook.
# Row r (new) receives the values that used to live in row $map[r] (old).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 3
    3  = 9
    4  = 8
    5  = 10
    6  = 16
    7  = 17
    8  = 6
    9  = 18
    10 = 24
    11 = 14
    12 = 11
    13 = 22
    14 = 2
    15 = 13
    16 = 4
    17 = 5
    18 = 26
    19 = 12
    20 = 21
    21 = 23
    22 = 25
    23 = 7
    24 = 15
    25 = 19
    26 = 20
}

$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# Snapshot all the original values for the affected columns/rows before
# writing anything back, since the permutation rewrites cells in place.
$orig = @{}
for ($r = 2; $r -le 26; $r++) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

foreach ($r in $map.Keys) {
    $src = $map[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $orig[$src][$c]
    }
}
